$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Shortname" column (E) with the navigators' first names.
# Fill data rows first (matches the order new shared strings appear in the
# saved file), then go back and fill in the header + first data row.
$ws.Range("E3").Value = "Mikhaela"
$ws.Range("E4").Value = "Chioma"
$ws.Range("E5").Value = "Evan"
$ws.Range("E6").Value = "Clare"
$ws.Range("E7").Value = "Caroline"
$ws.Range("E9").Value = "Shea"
$ws.Range("E8").Value = "Pat"
$ws.Range("E1").Value = "Shortname"
$ws.Range("E2").Value = "Meg"

# The D-column cells that previously carried a leftover/unused alignment
# style no longer need it - clear the formatting so the style goes away.
$ws.Range("D2").ClearFormats()
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D9").ClearFormats()

# Update the selected / active cell on the sheet (also drops the
# previous topLeftCell scroll position).
[void]$ws.Range("J7").Select()
